$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Ccl2"
$ws.Cells.Item(2, 3).Value = "Cxcr3"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 36.856666
$ws.Cells.Item(2, 8).Value = 110.569998
$ws.Cells.Item(2, 9).Value = 0.1389359711410444
$ws.Cells.Item(2, 10).Value = 0.1389359711410443
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 1.863797
$ws.Cells.Item(2, 14).Value = 5.591391
$ws.Cells.Item(2, 15).Value = 0.5455000708290748
$ws.Cells.Item(2, 16).Value = 0.5455000708290748
$ws.Cells.Item(2, 17).Value = 68.69334352080199
$ws.Cells.Item(2, 18).Value = 618.240091687218
$ws.Cells.Item(2, 19).Value = 0.075789582098146
$ws.Cells.Item(2, 20).Value = 0.07578958209814597

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Ccl2"
$ws.Cells.Item(3, 3).Value = "Cxcr3"
$ws.Cells.Item(3, 4).Value = "M2"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 36.856666
$ws.Cells.Item(3, 8).Value = 110.569998
$ws.Cells.Item(3, 9).Value = 0.1389359711410444
$ws.Cells.Item(3, 10).Value = 0.1389359711410443
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.552879
$ws.Cells.Item(3, 14).Value = 4.658637000000001
$ws.Cells.Item(3, 15).Value = 0.4544999291709252
$ws.Cells.Item(3, 16).Value = 0.4544999291709252
$ws.Cells.Item(3, 17).Value = 57.233942641414
$ws.Cells.Item(3, 18).Value = 515.1054837727261
$ws.Cells.Item(3, 19).Value = 0.06314638904289838
$ws.Cells.Item(3, 20).Value = 0.06314638904289836

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Ccl2"
$ws.Cells.Item(4, 3).Value = "Cxcr3"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 23.77965533333333
$ws.Cells.Item(4, 8).Value = 71.338966
$ws.Cells.Item(4, 9).Value = 0.089640487480229
$ws.Cells.Item(4, 10).Value = 0.08964048748022899
$ws.Cells.Item(4, 11).Value = 1
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 1.863797
$ws.Cells.Item(4, 14).Value = 5.591391
$ws.Cells.Item(4, 15).Value = 0.5455000708290748
$ws.Cells.Item(4, 16).Value = 0.5455000708290748
$ws.Cells.Item(4, 17).Value = 44.32045027130066
$ws.Cells.Item(4, 18).Value = 398.884052441706
$ws.Cells.Item(4, 19).Value = 0.04889889226961772
$ws.Cells.Item(4, 20).Value = 0.04889889226961771

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Ccl2"
$ws.Cells.Item(5, 3).Value = "Cxcr3"
$ws.Cells.Item(5, 4).Value = "M2"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 23.77965533333333
$ws.Cells.Item(5, 8).Value = 71.338966
$ws.Cells.Item(5, 9).Value = 0.089640487480229
$ws.Cells.Item(5, 10).Value = 0.08964048748022899
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.552879
$ws.Cells.Item(5, 14).Value = 4.658637000000001
$ws.Cells.Item(5, 15).Value = 0.4544999291709252
$ws.Cells.Item(5, 16).Value = 0.4544999291709252
$ws.Cells.Item(5, 17).Value = 36.92692739437134
$ws.Cells.Item(5, 18).Value = 332.342346549342
$ws.Cells.Item(5, 19).Value = 0.04074159521061129
$ws.Cells.Item(5, 20).Value = 0.04074159521061128

# Row 6
$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "Ccl2"
$ws.Cells.Item(6, 3).Value = "Cxcr3"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 199.1856246666667
$ws.Cells.Item(6, 8).Value = 597.556874
$ws.Cells.Item(6, 9).Value = 0.7508559835661449
$ws.Cells.Item(6, 10).Value = 0.7508559835661449
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 1.863797
$ws.Cells.Item(6, 14).Value = 5.591391
$ws.Cells.Item(6, 15).Value = 0.5455000708290748
$ws.Cells.Item(6, 16).Value = 0.5455000708290748
$ws.Cells.Item(6, 17).Value = 371.2415696968593
$ws.Cells.Item(6, 18).Value = 3341.174127271734
$ws.Cells.Item(6, 19).Value = 0.4095919922177667
$ws.Cells.Item(6, 20).Value = 0.4095919922177667

# Row 7
$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "Ccl2"
$ws.Cells.Item(7, 3).Value = "Cxcr3"
$ws.Cells.Item(7, 4).Value = "M2"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 199.1856246666667
$ws.Cells.Item(7, 8).Value = 597.556874
$ws.Cells.Item(7, 9).Value = 0.7508559835661449
$ws.Cells.Item(7, 10).Value = 0.7508559835661449
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.552879
$ws.Cells.Item(7, 14).Value = 4.658637000000001
$ws.Cells.Item(7, 15).Value = 0.4544999291709252
$ws.Cells.Item(7, 16).Value = 0.4544999291709252
$ws.Cells.Item(7, 17).Value = 309.3111736467487
$ws.Cells.Item(7, 18).Value = 2783.800562820738
$ws.Cells.Item(7, 19).Value = 0.3412639913483783
$ws.Cells.Item(7, 20).Value = 0.3412639913483783

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Ccl2"
$ws.Cells.Item(8, 3).Value = "Cxcr3"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 5.456122
$ws.Cells.Item(8, 8).Value = 16.368366
$ws.Cells.Item(8, 9).Value = 0.02056755781258178
$ws.Cells.Item(8, 10).Value = 0.02056755781258178
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 1.863797
$ws.Cells.Item(8, 14).Value = 5.591391
$ws.Cells.Item(8, 15).Value = 0.5455000708290748
$ws.Cells.Item(8, 16).Value = 0.5455000708290748
$ws.Cells.Item(8, 17).Value = 10.169103815234
$ws.Cells.Item(8, 18).Value = 91.52193433710599
$ws.Cells.Item(8, 19).Value = 0.01121960424354445
$ws.Cells.Item(8, 20).Value = 0.01121960424354445

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Ccl2"
$ws.Cells.Item(9, 3).Value = "Cxcr3"
$ws.Cells.Item(9, 4).Value = "M2"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 5.456122
$ws.Cells.Item(9, 8).Value = 16.368366
$ws.Cells.Item(9, 9).Value = 0.02056755781258178
$ws.Cells.Item(9, 10).Value = 0.02056755781258178
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.552879
$ws.Cells.Item(9, 14).Value = 4.658637000000001
$ws.Cells.Item(9, 15).Value = 0.4544999291709252
$ws.Cells.Item(9, 16).Value = 0.4544999291709252
$ws.Cells.Item(9, 17).Value = 8.472697275238
$ws.Cells.Item(9, 18).Value = 76.254275477142
$ws.Cells.Item(9, 19).Value = 0.009347953569037329
$ws.Cells.Item(9, 20).Value = 0.009347953569037327

